$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("A1").Value = "Hour"
$ws.Range("B1").Value = "Sum Count Of Purchased Products"

# Update the selection to column B, active cell B1
$ws.Range("B1:B1048576").Select()
